$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.730.36'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.809.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.555'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.56'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.90%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0713'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0929'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.066.32'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.814.79'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.05'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.646'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.762.71'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.84'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '255.16'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0805'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.50%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.73%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.19'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.55'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.47%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.98%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.28%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.447.70'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.39%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0193'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.17%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.07'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '85.49'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.966'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.21%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.12'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +7.30%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.06'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.962.98'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0492'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.41%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.17'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +8.94%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.13'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.89%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0127'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +9.13%  '
